$wb = $excel.ActiveWorkbook

# --- Update status text -----------------------------------------------
# The handoff status moved from "Ready for handoff" to "In Translation"
# everywhere it is reported: the Overview roll-up (columns E/F) and the
# per-locale detail sheets (column C on "zh-cn" / "de-de").
$lookAtWhole = [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation", $lookAtWhole) | Out-Null
}

# --- Resize status columns to fit the new, shorter text ---------------
# Overview: columns E (zh-cn) and F (de-de) show this status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# zh-cn / de-de detail sheets: column C is the Status column.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C1").ColumnWidth = 12.5

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C1").ColumnWidth = 12.5
